$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The question text in C2 was reworded.
$ws.Range("C2").Value = "Millise tunde viimane plokk sinus tekitas?"

# Move the active selection, matching the post-edit cursor position.
$ws.Range("C11").Select()
